$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two "Deck: College Of Winterhold" description lines ---
$ws.Range("B93").Value = "Consists of most mage cards. Pack a handfull of clear weathers & frostbites, has more spellswords than warriors"
$ws.Range("B94").Value = "Artifacts: Staff Of Magnus, Sigil Stone (Summon a Powerful Daedra)"

# --- Warriors sub-table: remove the blank spacer row (old row 100), pulling
#     rows 101-109 up by one, and fill in the new "Tolfdir" card in what is
#     now row 101 ---
$ws.Range("A100:H100").ClearContents()
$ws.Range("A101").Value = 3
$ws.Range("B101").Value = "Tolfdir (Spellsword)"
$ws.Range("C101").Value = 8
$ws.Range("D101").Value = "master alteration"
$ws.Range("E101").Value = "yes"
$ws.Range("G101").Value = "Nord"

$ws.Range("A102").Value = 4
$ws.Range("A103").Value = 5
$ws.Range("A104").Value = 6
$ws.Range("A105").Value = 7
$ws.Range("A106").Value = 8
$ws.Range("A107").Value = 9
$ws.Range("A108").Value = 10
$ws.Range("A109").Value = 11

# --- Mages sub-table: remove the blank spacer row (old row 117), pulling
#     rows 118-124 up by one, and fill in the new mage cards ---
$ws.Range("A117:H117").ClearContents()

$ws.Range("A118").Value = 5
$ws.Range("B118").Value = "Brelyna Maryon"
$ws.Range("C118").Value = 4
$ws.Range("D118").Value = "fire mage"
$ws.Range("E118").Value = "yes"

$ws.Range("A119").Value = 6
$ws.Range("B119").Value = "J'Zargo"
$ws.Range("C119").Value = 2
$ws.Range("D119").Value = "fire mage"
$ws.Range("E119").Value = "yes"

$ws.Range("A120").Value = 7
$ws.Range("B120").Value = "Phinis Gestor"
$ws.Range("C120").Value = 9
$ws.Range("D120").Value = "master conjurer"
$ws.Range("E120").Value = "yes"
$ws.Range("G120").Value = "breton"
$ws.Range("H120").Value = "summon a high level daedric warrior"

$ws.Range("A121").Value = 8
$ws.Range("B121").Value = "Faralda"
$ws.Range("C121").Value = 8
$ws.Range("D121").Value = "master destruction"
$ws.Range("E121").Value = "yes"
$ws.Range("G121").Value = "altmer"
$ws.Range("H121").Value = "Scorch but most likely might kill herself?"

$ws.Range("A122").Value = 9
$ws.Range("A123").Value = 10

$ws.Range("A124").Value = 11
$ws.Range("B124").Value = "Savos Aren"
$ws.Range("C124").Value = 11
$ws.Range("D124").Value = "frost mage"
$ws.Range("F124").Value = "yes"

# --- Shadow sub-table: fill in the first card (Drevis) ---
$ws.Range("B129").Value = "Drevis"
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = "master illusionist"
$ws.Range("E129").Value = "yes"
$ws.Range("F129").Value = "yes"
$ws.Range("G129").Value = "dunmer"
$ws.Range("H129").Value = "Spy card. (Do whenm spy cards are implemented)"

# --- New trailing note row ---
$ws.Range("B140").Value = "no images for mirabelle, collete, onmund"

# --- Restore the view state (scroll position / active cell) ---
$ws.Range("H122").Select()
